$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("5 November 2015")

# --- New rows of meeting-notes data appended below the existing ones ---
# Write order matters only insofar as it reproduces the natural order the
# rows/cells were authored in (A before B, row by row), with C5 filled
# in afterwards, matching how the sheet was actually edited.
$ws.Range("A4").Value = "UMUR PIUTANG"
$ws.Range("B4").Value = "0 TO 30"

$ws.Range("A5").Value = "CASH IN"
$ws.Range("B5").Value = "MASIH NGEBUG"

$ws.Range("B6").Value = "NO TELP"
$ws.Range("A6").Value = "DAFTAR SPP"

$ws.Range("B7").Value = "FORMAT SURAT"
$ws.Range("A7").Value = "PPJB"

$ws.Range("A8").Value = "DAFTAR KUITANSI"
$ws.Range("B8").Value = "JUDUL PERBAIKI"

$ws.Range("A9").Value = "MANUAL BOOK"
$ws.Range("B9").Value = "SCREEN SHOOT"

$ws.Range("A10").Value = "HAK AKSES "
$ws.Range("B10").Value = "HAK AKSES PER MODUL"

$ws.Range("A11").Value = "COLLECTION"
$ws.Range("B11").Value = "BELUM DI CEK (PENAGIHAN , PERHITUNGAN DENDA, CREATE TAGIHAN)"

$ws.Range("A12").Value = "NAMA PEJABAT, PEJABAT PPJB"
$ws.Range("B12").Value = "UBAH DIDATABASE"

$ws.Range("A14").Value = "MENU KARTU PEMBELI"
$ws.Range("B14").Value = "HAPUS"

$ws.Range("B13").Value = "HAPUS"
$ws.Range("A13").Value = "MENU RENCANA REALISASI SELURUH BLOK"

$ws.Range("C13").Value = "OK"
$ws.Range("C14").Value = "OK"

$ws.Range("C5").Value = "OTW"

# --- Column widths widened (best-fit) to accommodate the new, longer text ---
$ws.Columns.Item(1).ColumnWidth = 38.15
$ws.Columns.Item(2).ColumnWidth = 63.6

# --- Selected cell in the sheet moved to C15 ---
$ws.Range("C15").Select() | Out-Null
